$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column I ("Tipo di supporto") contains values CERTIFICATO / TESSERINO in rows 2-10.
# Replace them with lowercase equivalents: certificato / tesserino.
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $val = $cell.Value2
    if ($val -eq "CERTIFICATO") {
        $cell.Value2 = "certificato"
    } elseif ($val -eq "TESSERINO") {
        $cell.Value2 = "tesserino"
    }
}

# Reflect the final selection location left by the editor.
[void]$ws.Range("I20").Select()
